# "update harian tusbung sendiri"
# Mark the currently-filtered ("JHON") rows on sheet "8" as paid ("lunas")
# in column P, wherever they still read "blm lunas".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8")

$ranges = @(
    "P254:P267","P269","P271","P273:P275","P277","P279","P282:P283",
    "P285:P290","P294:P302","P304:P306","P310","P312:P323","P325",
    "P327:P331","P333:P341","P343","P345:P352","P354:P361"
)

foreach ($addr in $ranges) {
    $ws.Range($addr).Value = "lunas"
}

# Reflect where the user ended up after the bulk edit: scrolled further
# down the filtered list, with the newly-updated P column selected.
$ws.Activate()
[void]$ws.Range("P254:P476").Select()
$excel.ActiveWindow.ScrollRow = 355
$excel.ActiveWindow.ScrollColumn = 1
